# Burndown chart and sprint log update
# The sprint backlog kanban table (ToDo / In Progress / Reviewing / Done) gets
# refreshed: several task cards move to a different status column, and a few
# task labels get split so the column only holds the assignee's name while the
# task description moves into its own status cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Folha1")
$ws.Activate()

# Row 2: "Preparing the workplace" moves from In Progress (B) to Done (D)
$ws.Range("B2").Value = ""
$ws.Range("D2").Value = "Preparing the workplace"

# Row 4: "Identifie 3 code smells Martim Costa" (ToDo) splits into the
# assignee name (ToDo) and the task, which now also shows up in Reviewing
$ws.Range("A4").Value = " Martim Costa"
$ws.Range("C4").Value = "Identifie 3 code smells"
$ws.Range("D4").Value = ""

# Row 7: "Identifie 3 code smells - Pedro Arruda" -> name stays in ToDo,
# task moves to Reviewing
$ws.Range("A7").Value = " Pedro Arruda"
$ws.Range("C7").Value = "Identifie 3 code smells"

# Row 8: "Identifie 3 code smells - Bernardo Carvalho" -> name stays in ToDo,
# task moves to Reviewing
$ws.Range("A8").Value = " Bernardo Carvalho"
$ws.Range("C8").Value = "Identifie 3 code smells"

# Row 10: "Identifie 3 design paterns" moves from Done (D) to Reviewing (C)
$ws.Range("C10").Value = "Identifie 3 design paterns"
$ws.Range("D10").Value = ""

# Row 13: "Identifie 3 design paterns - Pedro Arruda" -> name stays in ToDo,
# task moves to Reviewing
$ws.Range("A13").Value = "Pedro Arruda"
$ws.Range("C13").Value = "Identifie 3 design paterns"

# Row 14: "Identifie 3 design paterns - Bernardo Carvalho" -> name stays in
# ToDo, task moves to Reviewing
$ws.Range("A14").Value = "Bernardo Carvalho"
$ws.Range("C14").Value = "Identifie 3 design paterns "

# Cursor / selection ends on C14, matching the last edited cell
$ws.Range("C14").Select()

